# HighlightPoints.pptx - "Preserve object order when swapping" (#1302/#1333)
#
# The PPT Labs "Highlight Bullets" add-in keeps one background-highlight
# shape per bullet line plus a single text shape that carries the bullet
# text itself. When the add-in re-swaps/re-highlights points it must move
# the text shape backwards in z-order so it lines up with the matching
# background shape again. This script reproduces that z-order swap for
# the two affected slides (slide 14 and slide 17 of the deck).

$msoSendBackward = 3

$p = $ppt.ActivePresentation

# ---- Slide 14 ------------------------------------------------------------
# Text shape "HighlightBackgroundShape961d3f1e-..." must move from just
# after the 7th (last) PPTLabsHighlightBackgroundShape back to just before
# the 5th one - i.e. three steps backward in z-order.
$s14 = $p.Slides.Item(14)
$shape14 = $s14.Shapes.Item("HighlightBackgroundShape961d3f1e-187d-4b94-8f93-1be85193a105")
for ($i = 0; $i -lt 3; $i++) {
    $shape14.ZOrder($msoSendBackward)
}

# ---- Slide 17 ------------------------------------------------------------
# Text shape "HighlightBackgroundShapeea60d9b9-..." moves four steps
# backward (from just before the last background shape group to just
# before the 5th background shape).
$s17 = $p.Slides.Item(17)
$shapeA17 = $s17.Shapes.Item("HighlightBackgroundShapeea60d9b9-ed30-4da3-a1da-6902e7f48621")
for ($i = 0; $i -lt 4; $i++) {
    $shapeA17.ZOrder($msoSendBackward)
}

# Text shape "HighlightBackgroundShapec485a6a5-..." moves one step
# backward (swaps with the background shape that was already out of
# numeric order).
$shapeB17 = $s17.Shapes.Item("HighlightBackgroundShapec485a6a5-1d97-42dc-ab66-ce7355d7009d")
$shapeB17.ZOrder($msoSendBackward)
